$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Force the value to be written as plain text (not auto-converted into a
    # date serial number by Excel's input parser), while leaving the cell's
    # number format/style untouched afterwards.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# B2: cert expiry date moved out a week
Set-TextValue "B2" "2027-08-05"

# B3: cert expiry date moved out a week
Set-TextValue "B3" "2027-08-05"

# E5: expired day count incremented by one day
$ws.Range("E5").Value = "❌ EXPIRED 3769 days ago"

# B6: cert expiry date pushed out roughly a year; no longer "expiring soon"
Set-TextValue "B6" "2026-07-17"
$ws.Range("E6").Value = ""

# B8: cert expiry date pushed out roughly a year; no longer "expiring soon"
Set-TextValue "B8" "2026-07-17"
$ws.Range("E8").Value = ""

# E9: days-until-expiry count decremented by one day
$ws.Range("E9").Value = "⚠️ Expires in 29 days"
